$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the VIN value (column A) used by rows 2-5 (shared string) to the new valid VIN
$ws.Range("A2").Value = "4T1BE30K&6"
$ws.Range("A3").Value = "4T1BE30K&6"
$ws.Range("A4").Value = "4T1BE30K&6"
$ws.Range("A5").Value = "4T1BE30K&6"

# Update the MAKE_TEXT value in E2 from TOYOTA to TOYOTA_UPDATED
$ws.Range("E2").Value = "TOYOTA_UPDATED"

# Reset the selection to E3
$ws.Range("E3").Select() | Out-Null
